# CORE_holdings.xlsx update
# - Bump the "Model holdings provided as of" disclosure date from 2021-05-14 to 2021-05-17
# - Refresh the Weight (D) and Percent Change (E) figures for rows 2-8 with the latest values
#
# The worksheet ships with sheet-level protection (no cells are explicitly unlocked), so the
# protected cells must be unlocked for editing and then the protection state restored.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$wasProtected = $ws.ProtectContents
if ($wasProtected) {
    $ws.Unprotect()
}

# --- Update the confidential disclosure date in the footnote text (A11) ---
$ws.Range("A11").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-17 for illustrative purposes only and are subject to change."

# --- Refresh Weight (D) and Percent Change (E) values for rows 2-8 ---
$ws.Range("D2").Value = 0.5014986006784071
$ws.Range("E2").Value = -0.0006674676278199954

$ws.Range("D3").Value = 0.2427453519985622
$ws.Range("E3").Value = -0.005155398438650649

$ws.Range("D4").Value = 0.09518467906894477
$ws.Range("E4").Value = -0.004657603222558038

$ws.Range("D5").Value = 0.1036439418361496
$ws.Range("E5").Value = 0.003394806863014965

$ws.Range("D6").Value = 0.03007756891765827
$ws.Range("E6").Value = 0.005501280470454395

$ws.Range("D7").Value = 0.02684985750027815
$ws.Range("E7").Value = -0.003817983481377629

$ws.Range("E8").Value = -0.001614711562875071

# Restore sheet protection to match the original state
if ($wasProtected) {
    $ws.Protect()
}
